# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 489
$wsExhibit.Range("F3").Value = 5810
$wsExhibit.Range("F4").Value = 391
$wsExhibit.Range("F5").Value = 76
$wsExhibit.Range("F9").Value = 547
$wsExhibit.Range("F10").Value = 25

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 489
$wsAll.Range("F3").Value = 5810
$wsAll.Range("F4").Value = 391
$wsAll.Range("F6").Value = 76
$wsAll.Range("F11").Value = 547
$wsAll.Range("F12").Value = 25

$wb.Save()
